# v1.5 closed publish and upload reviews
#
# 1) On REVIEW-SHEET: mark the Publish/Upload related reviewer-verification
#    cells (J5:J10 -> SRS-PUB-001..006) as "closed" (they were "open").
# 2) On VERSION-HISTORY: append a new history row (row 8) documenting the
#    v1.6 update ("Close reviewer verification for Publish and upload
#    feature", owner Gehad Ashry, date 2025-04-18).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("REVIEW-SHEET")
$ws2 = $wb.Worksheets.Item("VERSION-HISTORY")

# --- 1. Close the reviewer verification for the publish/upload rows ---
$ws1.Range("J5").Value  = "closed"
$ws1.Range("J6").Value  = "closed"
$ws1.Range("J7").Value  = "closed"
$ws1.Range("J8").Value  = "closed"
$ws1.Range("J9").Value  = "closed"
$ws1.Range("J10").Value = "closed"

# --- 2. Append the new version-history row, copying formatting from the ---
#        previous (v1.5) row so fonts/fills/borders/number-format match.
$ws2.Range("A7:D7").Copy()
$ws2.Range("A8:D8").PasteSpecial(-4122)

$ws2.Range("A8").Value = "v1.6"
$ws2.Range("B8").Value = "Gehad Ashry"
$ws2.Range("C8").Value = "Close reviewer verification for Publish and upload feature"
$ws2.Range("D8").Value = 45765
$ws2.Rows.Item(8).RowHeight = 18.75

# --- 3. Leave the selection/active sheet the way the author left it: ---
#        REVIEW-SHEET active, last cell touched (J10) selected, and the
#        version-history sheet remembers its own last selection too.
$ws2.Range("D11").Select()
$ws1.Select()
$ws1.Range("J10").Select()
